# Update the TPM-derived values in the "Fgf10-Fgfr2" sheet.
# The underlying TPM values were re-run and the receptor-expression /
# specificity-derived columns (K-T) were recalculated for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster = ECs)
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.05646129082
$ws.Range("R2").Value = 0.50815161738
$ws.Range("S2").Value = 0.02099032928903418
$ws.Range("T2").Value = 0.02099032928903418

# Row 3 (Target cluster = FAPs)
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("S3").Value = 0.5358731102718634
$ws.Range("T3").Value = 0.5358731102718634

# Row 4 (Target cluster = MuSCs)
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("S4").Value = 0.4431365604391025
$ws.Range("T4").Value = 0.4431365604391026
